# Applies the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list on Fri May 31 16:30:36 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.848.80"
$ws.Range("E2").Value = "  -2.30%  "

$ws.Range("D3").Value = "3.737.46"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.30%  "

$ws.Range("D7").Value = "3.739.30"
$ws.Range("E7").Value = "  -0.94%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.66%  "

$ws.Range("D15").Value = "4.359.90"
$ws.Range("E15").Value = "  -1.07%  "

$ws.Range("D16").Value = "3.748.84"
$ws.Range("E16").Value = "  -0.81%  "

$ws.Range("D17").Value = "66.853.83"
$ws.Range("E17").Value = "  -2.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.62%  "

$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "451.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.688"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000146"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.42%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.90%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.26%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("E29").Value = "  -2.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.54%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.24%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.81%  "

$ws.Range("D36").Value = "3.688.57"
$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0983"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.985"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.24%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.294"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.85%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "146.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "383.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.92%  "

$ws.Range("D51").Value = "2.726.91"
$ws.Range("E51").Value = "  +1.36%  "
